# "Changed few mistakes in presentation"
#
# Slide 1, subtitle placeholder ("Подзаголовок 2"): the word
# "последованности" was misspelled; correct it to "последованости".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$bad = $tr.Find("последованности")
$bad.Text = "последованости"
